$d = $word.ActiveDocument

# The cover-page placeholder paragraph currently reads "<Tag der Promotion>"
# (three runs: "<", "Tag der Promotion", ">"). The commit removes that text
# entirely, leaving the paragraph (and its paragraph-mark run properties)
# in place but empty.
$found = $d.Content.Find.Execute("<Tag der Promotion>", $true, $true, $false, $false, $false,
                                  $true, 1, $false, "", 2)

if (-not $found) {
    throw "Could not find the '<Tag der Promotion>' placeholder text to remove."
}

Write-Host "Removed '<Tag der Promotion>' placeholder: $found"
